# Simulated Wild Card round and logged it
# Appends the new game's play-by-play / per-kick logs to the running logs
# (YDS and ST sheets) and updates the aggregated season totals on the
# OFF, DEF, ST, TURNS and PEN sheets to include the new game.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: append this game's rushing/passing yardage-per-play logs
# ---------------------------------------------------------------------
$wsYDS = $wb.Worksheets.Item("YDS")

$rushOff = $wsYDS.Cells.Item(2,2).Value2
$wsYDS.Cells.Item(2,2).Value = $rushOff + " 6 1 0 3 1 -2 -1 -1 1 7 6 5 11 2 6 3 -3 15 3"

$rushDef = $wsYDS.Cells.Item(2,3).Value2
$wsYDS.Cells.Item(2,3).Value = $rushDef + " -3 5 3 5 2 4 1 3 16 1 3 16 13 -1 2 8 -1 5 8 8 0 2 3 3 1 8 -2 10 2 0 11"

$passOff = $wsYDS.Cells.Item(3,2).Value2
$wsYDS.Cells.Item(3,2).Value = $passOff + " 7 8 14 6 -3 25 4 7 9 -5 20 28 24 8 13 37 30 7 5 5 21 5 26"

$passDef = $wsYDS.Cells.Item(3,3).Value2
$wsYDS.Cells.Item(3,3).Value = $passDef + " 5 2 29 13 4 9 19 17 9 15 0 5 14 24 14 1 14 4 10 36 24 3 16"

# ---------------------------------------------------------------------
# OFF sheet: roll the new game's down & distance splits into the totals
# ---------------------------------------------------------------------
$wsOFF = $wb.Worksheets.Item("OFF")

$wsOFF.Cells.Item(2,3).Value  = 316
$wsOFF.Cells.Item(2,4).Value  = 14
$wsOFF.Cells.Item(2,5).Value  = 23
$wsOFF.Cells.Item(2,6).Value  = 134
$wsOFF.Cells.Item(2,7).Value  = 96
$wsOFF.Cells.Item(2,8).Value  = 9
$wsOFF.Cells.Item(2,9).Value  = 13
$wsOFF.Cells.Item(2,10).Value = 66
$wsOFF.Cells.Item(2,12).Value = 584
$wsOFF.Cells.Item(2,13).Value = 404
$wsOFF.Cells.Item(2,15).Value = 37
$wsOFF.Cells.Item(2,16).Value = 18
$wsOFF.Cells.Item(2,17).Value = 1025

$wsOFF.Cells.Item(3,2).Value  = 19
$wsOFF.Cells.Item(3,3).Value  = 348
$wsOFF.Cells.Item(3,5).Value  = 66
$wsOFF.Cells.Item(3,6).Value  = 180
$wsOFF.Cells.Item(3,7).Value  = 62
$wsOFF.Cells.Item(3,8).Value  = 63
$wsOFF.Cells.Item(3,9).Value  = 118
$wsOFF.Cells.Item(3,10).Value = 120
$wsOFF.Cells.Item(3,14).Value = 48

# ---------------------------------------------------------------------
# DEF sheet: roll the new game's down & distance splits into the totals
# ---------------------------------------------------------------------
$wsDEF = $wb.Worksheets.Item("DEF")

$wsDEF.Cells.Item(2,2).Value  = 8
$wsDEF.Cells.Item(2,3).Value  = 465
$wsDEF.Cells.Item(2,5).Value  = 24
$wsDEF.Cells.Item(2,6).Value  = 124
$wsDEF.Cells.Item(2,7).Value  = 127
$wsDEF.Cells.Item(2,10).Value = 66
$wsDEF.Cells.Item(2,12).Value = 556
$wsDEF.Cells.Item(2,13).Value = 379
$wsDEF.Cells.Item(2,17).Value = 1098

$wsDEF.Cells.Item(3,3).Value  = 335
$wsDEF.Cells.Item(3,4).Value  = 12
$wsDEF.Cells.Item(3,5).Value  = 57
$wsDEF.Cells.Item(3,6).Value  = 207
$wsDEF.Cells.Item(3,7).Value  = 68
$wsDEF.Cells.Item(3,8).Value  = 54
$wsDEF.Cells.Item(3,9).Value  = 109
$wsDEF.Cells.Item(3,10).Value = 103
$wsDEF.Cells.Item(3,14).Value = 30

# ---------------------------------------------------------------------
# ST sheet: roll the new game's totals in and append the per-kick logs
# ---------------------------------------------------------------------
$wsST = $wb.Worksheets.Item("ST")

$wsST.Cells.Item(2,2).Value  = 147
$wsST.Cells.Item(2,4).Value  = 141
$wsST.Cells.Item(2,6).Value  = 185
$wsST.Cells.Item(2,7).Value  = 167
$wsST.Cells.Item(2,8).Value  = 7
$wsST.Cells.Item(2,9).Value  = 4
$wsST.Cells.Item(2,10).Value = 76
$wsST.Cells.Item(2,11).Value = 73

$wsST.Cells.Item(3,2).Value  = 84

$koD = $wsST.Cells.Item(4,2).Value2
$wsST.Cells.Item(4,2).Value = $koD + " 63 63 49"

$koRA = $wsST.Cells.Item(5,2).Value2
$wsST.Cells.Item(5,2).Value = $koRA + " 25 21 12"

$koRM = $wsST.Cells.Item(6,2).Value2
$wsST.Cells.Item(6,2).Value = $koRM + " 16 24 18"

$ptTB = $wsST.Cells.Item(3,4).Value2
$wsST.Cells.Item(3,4).Value = $ptTB + " 39 33 42 50 47"

$ptD = $wsST.Cells.Item(4,4).Value2
$wsST.Cells.Item(4,4).Value = $ptD + " 1 0 24 8 10"

$ptRA = $wsST.Cells.Item(5,4).Value2
$wsST.Cells.Item(5,4).Value = $ptRA + " 13 3 8 0"

# ---------------------------------------------------------------------
# TURNS sheet: add this game's turnovers
# ---------------------------------------------------------------------
$wsTURNS = $wb.Worksheets.Item("TURNS")

$wsTURNS.Cells.Item(2,5).Value = 13
$wsTURNS.Cells.Item(3,4).Value = 17

# ---------------------------------------------------------------------
# PEN sheet: add this game's penalties
# ---------------------------------------------------------------------
$wsPEN = $wb.Worksheets.Item("PEN")

$wsPEN.Cells.Item(2,2).Value = 34
$wsPEN.Cells.Item(2,4).Value = 17
$wsPEN.Cells.Item(4,4).Value = 14
